$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (pair FAPs/MuSCs,Csf2,Sdc2,Neutrophils removed)
$ws.Rows(11).Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Csf2"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2300786666666667
$ws.Range("H2").Value = 0.690236
$ws.Range("I2").Value = 0.2137022699341201
$ws.Range("J2").Value = 0.2304051507653011
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.020434
$ws.Range("N2").Value = 8.040868
$ws.Range("O2").Value = 0.06420203970218387
$ws.Range("P2").Value = 0.04878916602310146
$ws.Range("Q2").Value = 0.9250160941413332
$ws.Range("R2").Value = 5.550096564847999
$ws.Range("S2").Value = 0.01372012161875719
$ws.Range("T2").Value = 0.011241275153266

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Csf2"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2300786666666667
$ws.Range("H3").Value = 0.690236
$ws.Range("I3").Value = 0.2137022699341201
$ws.Range("J3").Value = 0.2304051507653011
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 39.565288
$ws.Range("N3").Value = 118.695864
$ws.Range("O3").Value = 0.6318154186847339
$ws.Range("P3").Value = 0.7202048603398876
$ws.Range("Q3").Value = 9.103128709322666
$ws.Range("R3").Value = 81.92815838390399
$ws.Range("S3").Value = 0.1350203891523041
$ws.Range("T3").Value = 0.1659389094285144

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Csf2"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2300786666666667
$ws.Range("H4").Value = 0.690236
$ws.Range("I4").Value = 0.2137022699341201
$ws.Range("J4").Value = 0.2304051507653011
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.035871
$ws.Range("N4").Value = 38.071742
$ws.Range("O4").Value = 0.3039825416130822
$ws.Range("P4").Value = 0.2310059736370109
$ws.Range("Q4").Value = 4.379747818518666
$ws.Range("R4").Value = 26.278486911112
$ws.Range("S4").Value = 0.06496175916305881
$ws.Range("T4").Value = 0.05322496618352068

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Csf2"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.2341465
$ws.Range("H5").Value = 0.468293
$ws.Range("I5").Value = 0.2174805655477089
$ws.Range("J5").Value = 0.1563191709318771
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.020434
$ws.Range("N5").Value = 8.040868
$ws.Range("O5").Value = 0.06420203970218387
$ws.Range("P5").Value = 0.04878916602310146
$ws.Range("Q5").Value = 0.941370549581
$ws.Range("R5").Value = 3.765482198324
$ws.Range("S5").Value = 0.01396269590374741
$ws.Range("T5").Value = 0.007626681983188929

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Csf2"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.2341465
$ws.Range("H6").Value = 0.468293
$ws.Range("I6").Value = 0.2174805655477089
$ws.Range("J6").Value = 0.1563191709318771
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 39.565288
$ws.Range("N6").Value = 118.695864
$ws.Range("O6").Value = 0.6318154186847339
$ws.Range("P6").Value = 0.7202048603398876
$ws.Range("Q6").Value = 9.264073706692001
$ws.Range("R6").Value = 55.584442240152
$ws.Range("S6").Value = 0.1374075745773184
$ws.Range("T6").Value = 0.1125818266694396

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Csf2"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.2341465
$ws.Range("H7").Value = 0.468293
$ws.Range("I7").Value = 0.2174805655477089
$ws.Range("J7").Value = 0.1563191709318771
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.035871
$ws.Range("N7").Value = 38.071742
$ws.Range("O7").Value = 0.3039825416130822
$ws.Range("P7").Value = 0.2310059736370109
$ws.Range("Q7").Value = 4.457182569101501
$ws.Range("R7").Value = 17.828730276406
$ws.Range("S7").Value = 0.06611029506664308
$ws.Range("T7").Value = 0.03611066227924862

# Row 8
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Csf2"
$ws.Range("C8").Value = "Sdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6124066666666667
$ws.Range("H8").Value = 1.83722
$ws.Range("I8").Value = 0.5688171645181709
$ws.Range("J8").Value = 0.6132756783028217
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.020434
$ws.Range("N8").Value = 8.040868
$ws.Range("O8").Value = 0.06420203970218387
$ws.Range("P8").Value = 0.04878916602310146
$ws.Range("Q8").Value = 2.462140584493333
$ws.Range("R8").Value = 14.77284350696
$ws.Range("S8").Value = 0.03651922217967926
$ws.Range("T8").Value = 0.02992120888664653

# Row 9
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Csf2"
$ws.Range("C9").Value = "Sdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6124066666666667
$ws.Range("H9").Value = 1.83722
$ws.Range("I9").Value = 0.5688171645181709
$ws.Range("J9").Value = 0.6132756783028217
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 39.565288
$ws.Range("N9").Value = 118.695864
$ws.Range("O9").Value = 0.6318154186847339
$ws.Range("P9").Value = 0.7202048603398876
$ws.Range("Q9").Value = 24.23004613978667
$ws.Range("R9").Value = 218.07041525808
$ws.Range("S9").Value = 0.3593874549551113
$ws.Range("T9").Value = 0.4416841242419335

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Csf2"
$ws.Range("C10").Value = "Sdc2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6124066666666667
$ws.Range("H10").Value = 1.83722
$ws.Range("I10").Value = 0.5688171645181709
$ws.Range("J10").Value = 0.6132756783028217
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.035871
$ws.Range("N10").Value = 38.071742
$ws.Range("O10").Value = 0.3039825416130822
$ws.Range("P10").Value = 0.2310059736370109
$ws.Range("Q10").Value = 11.65769430620667
$ws.Range("R10").Value = 69.94616583723999
$ws.Range("S10").Value = 0.1729104873833803
$ws.Range("T10").Value = 0.1416703451742416
